$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("generator_file")

# Row 7: duplicate generator (plant 55350) correction - prime mover CT
$ws.Cells.Item(7, 1).Value = 55350
$ws.Cells.Item(7, 2).NumberFormat = "@"
$ws.Cells.Item(7, 2).Value = "1"
$ws.Cells.Item(7, 3).Value = "prime_mover"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "CT"

# Row 8: duplicate generator (plant 55350) correction - prime mover CA
$ws.Cells.Item(8, 1).Value = 55350
$ws.Cells.Item(8, 2).NumberFormat = "@"
$ws.Cells.Item(8, 2).Value = "3"
$ws.Cells.Item(8, 3).Value = "prime_mover"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "CA"

$ws.Range("B9").Select()
